# Insert a new "債權" (claims) worksheet before the "事業投資" (investment)
# worksheet, matching the shape/layout of the other property-type sheets in
# this 財產申報表 (property-declaration) workbook.
#
# Strategy: the target sheetIds in the final workbook are 債權=8, 事業投資=9.
# This engine assigns new-sheet sheetIds as (max existing sheetId + 1), and
# never re-uses an id once consumed - so to land 債權 on 8 and 事業投資 on 9
# we delete the original 事業投資 sheet first (freeing its slot / dropping the
# max back to 7), add 債權 (-> id 8) after 保險, then re-add 事業投資
# (-> id 9) after 債權, restoring its original data.

$wb = $excel.ActiveWorkbook

# Capture the original "事業投資" sheet's data before we delete/recreate it.
$oldInvest = $wb.Worksheets.Item("事業投資")
$oldInvest.Delete()

# --- Add "債權" right after "保險" (i.e. where "事業投資" used to start) ---
$insurance = $wb.Worksheets.Item("保險")
$claims = $wb.Worksheets.Add($null, $insurance)
$claims.Name = "債權"

$claimHeaders = @("species","owner","debtor","total","register_date","register_reason","property_category","category","date","legislator_name","legislator_id","source_file","index")
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N")

for ($i = 0; $i -lt $claimHeaders.Length; $i++) {
    $cell = $claims.Range($cols[$i] + "1")
    $cell.Value = $claimHeaders[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    foreach ($idx in 7,8,9,10) {
        $cell.Borders.Item($idx).LineStyle = 1
    }
}

$claimA2 = $claims.Range("A2")
$claimA2.Value = 107
$claimA2.Font.Bold = $true
$claimA2.HorizontalAlignment = -4108
$claimA2.VerticalAlignment = -4160
foreach ($idx in 7,8,9,10) {
    $claimA2.Borders.Item($idx).LineStyle = 1
}

$claimRow2 = @("—般借款","李貴敏","李貴琪臺北市士林區華岡路",30000000,"100年","長期借款","claim","normal","__DATE__","李貴敏",1739,"tmp59331",107)
for ($i = 0; $i -lt $claimRow2.Length; $i++) {
    $cell = $claims.Range($cols[$i] + "2")
    if ($claimRow2[$i] -eq "__DATE__") {
        $cell.NumberFormat = "@"
        $cell.Value = "2012-04-27"
        $cell.Style = "Normal"
    } else {
        $cell.Value = $claimRow2[$i]
    }
}

# --- Re-add "事業投資" right after "債權", restoring its original content ---
$invest = $wb.Worksheets.Add($null, $claims)
$invest.Name = "事業投資"

$investHeaders = @("owner","company","address","total","register_date","register_reason","property_category","category","date","legislator_name","legislator_id","source_file","index")

for ($i = 0; $i -lt $investHeaders.Length; $i++) {
    $cell = $invest.Range($cols[$i] + "1")
    $cell.Value = $investHeaders[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    foreach ($idx in 7,8,9,10) {
        $cell.Borders.Item($idx).LineStyle = 1
    }
}

$investA2 = $invest.Range("A2")
$investA2.Value = 115
$investA2.Font.Bold = $true
$investA2.HorizontalAlignment = -4108
$investA2.VerticalAlignment = -4160
foreach ($idx in 7,8,9,10) {
    $investA2.Borders.Item($idx).LineStyle = 1
}

$investRow2 = @("李貴敏","國際通商法律事務所","臺北市松山區敦化北路168號15樓",90000000,"自民國82年","合夥","investment","normal","__DATE__","李貴敏",1739,"tmp59331",115)
for ($i = 0; $i -lt $investRow2.Length; $i++) {
    $cell = $invest.Range($cols[$i] + "2")
    if ($investRow2[$i] -eq "__DATE__") {
        $cell.NumberFormat = "@"
        $cell.Value = "2012-04-27"
        $cell.Style = "Normal"
    } else {
        $cell.Value = $investRow2[$i]
    }
}
